$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8496.25
$ws.Range("I62").Value = 2002.5
$ws.Range("J62").Value = 14990
$ws.Range("K62").Value = 2002.5
$ws.Range("L62").Value = 14990
$ws.Range("M62").Value = -1378.5
$ws.Range("N62").Value = -16238

$ws.Range("H65").Value = 8496.25
$ws.Range("I65").Value = 2002.5
$ws.Range("J65").Value = 14990
$ws.Range("K65").Value = 10012.5
$ws.Range("L65").Value = 74950
$ws.Range("M65").Value = -6892.5
$ws.Range("N65").Value = -81190

$ws.Range("H100").Value = 2186.3125
$ws.Range("I100").Value = 1617.5
$ws.Range("K100").Value = 1617.5
$ws.Range("M100").Value = -1076.5

$ws.Range("H116").Value = 19850.834
$ws.Range("I116").Value = 52252.5
$ws.Range("J116").Value = 3650
$ws.Range("K116").Value = 52252.5
$ws.Range("L116").Value = 3650
$ws.Range("M116").Value = -48810.5
$ws.Range("N116").Value = -10534

$ws.Range("H127").Value = 1496.4546
$ws.Range("I127").Value = 1496.4546
$ws.Range("K127").Value = 4489.3638
$ws.Range("M127").Value = 470.6361999999999

$ws.Range("H138").Value = 3155.1904
$ws.Range("J138").Value = 2393.926
$ws.Range("L138").Value = 7181.778
$ws.Range("N138").Value = -17461.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3752838
$ws.Range("I45").Value = 9002671
$ws.Range("J45").Value = 2957.1428
$ws.Range("K45").Value = 9002671
$ws.Range("L45").Value = 2957.1428
$ws.Range("M45").Value = -9002294
$ws.Range("N45").Value = -3711.1428

$ws.Range("H97").Value = 1307.9286
$ws.Range("I97").Value = 1246.1
$ws.Range("K97").Value = 1246.1
$ws.Range("M97").Value = -750.0999999999999

$ws.Range("H132").Value = 1429.4546
$ws.Range("I132").Value = 1138.4103
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 3415.2309
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -885.2309
$ws.Range("N132").Value = -16158.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H86").Value = 2262.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2262.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2262.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4508.5

$ws.Range("H89").Value = 2262.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2262.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 11312.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -22544.5

$ws.Range("H107").Value = 3647.4285
$ws.Range("I107").Value = 3547.4167
$ws.Range("J107").Value = 4247.5
$ws.Range("K107").Value = 3547.4167
$ws.Range("L107").Value = 4247.5
$ws.Range("M107").Value = -1627.4167
$ws.Range("N107").Value = -8087.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1384.2881
$ws.Range("I31").Value = 692.913
$ws.Range("K31").Value = 692.913
$ws.Range("M31").Value = -397.913

$ws.Range("H34").Value = 1384.2881
$ws.Range("I34").Value = 692.913
$ws.Range("K34").Value = 692.913
$ws.Range("M34").Value = -490.913

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H132").Value = 2331.0264
$ws.Range("I132").Value = 1708.7931
$ws.Range("K132").Value = 5126.379300000001
$ws.Range("M132").Value = -2596.379300000001

$ws.Range("H134").Value = 1552.78
$ws.Range("I134").Value = 861.4737
$ws.Range("K134").Value = 2584.4211
$ws.Range("M134").Value = -49.42110000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2333
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2333
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7919

$ws.Range("H68").Value = 2870.2917
$ws.Range("J68").Value = 2908.0435
$ws.Range("L68").Value = 8724.130500000001
$ws.Range("N68").Value = -10346.1305

$ws.Range("H71").Value = 2870.2917
$ws.Range("J71").Value = 2908.0435
$ws.Range("L71").Value = 26172.3915
$ws.Range("N71").Value = -34284.3915

$ws.Range("H107").Value = 1777.409
$ws.Range("J107").Value = 1915.9375
$ws.Range("L107").Value = 5747.8125
$ws.Range("N107").Value = -9587.8125

$ws.Range("H113").Value = 1470.3077
$ws.Range("I113").Value = 3780.6667
$ws.Range("J113").Value = 777.2
$ws.Range("K113").Value = 11342.0001
$ws.Range("L113").Value = 2331.6
$ws.Range("M113").Value = -9172.000100000001
$ws.Range("N113").Value = -6671.6

$ws.Range("H131").Value = 10885861
$ws.Range("J131").Value = 17401.883
$ws.Range("L131").Value = 52205.649
$ws.Range("N131").Value = -62285.649

$ws.Range("H141").Value = 1950.4546
$ws.Range("I141").Value = 1950.4546
$ws.Range("K141").Value = 5851.3638
$ws.Range("M141").Value = -671.3638000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3085230.8

$ws.Range("H70").Value = 5040
$ws.Range("I70").Value = 5400
$ws.Range("K70").Value = 5400
$ws.Range("M70").Value = -5130

$ws.Range("H73").Value = 5040
$ws.Range("I73").Value = 5400
$ws.Range("K73").Value = 5400
$ws.Range("M73").Value = -4464

$ws.Range("H80").Value = 1989
$ws.Range("I80").Value = 1986.4
$ws.Range("J80").Value = 1993.3334
$ws.Range("K80").Value = 1986.4
$ws.Range("L80").Value = 1993.3334
$ws.Range("M80").Value = -988.4000000000001
$ws.Range("N80").Value = -3989.3334

$ws.Range("H83").Value = 1989
$ws.Range("I83").Value = 1986.4
$ws.Range("J83").Value = 1993.3334
$ws.Range("K83").Value = 9932
$ws.Range("L83").Value = 9966.666999999999
$ws.Range("M83").Value = -4940
$ws.Range("N83").Value = -19950.667

$ws.Range("H97").Value = 1076.4242
$ws.Range("I97").Value = 809.5417
$ws.Range("K97").Value = 809.5417
$ws.Range("M97").Value = -313.5417

$ws.Range("H102").Value = 2877.9333
$ws.Range("I102").Value = 3055.8
$ws.Range("K102").Value = 3055.8
$ws.Range("M102").Value = -1433.8

$ws.Range("H132").Value = 1750938.6
$ws.Range("I132").Value = 2138480.5
$ws.Range("K132").Value = 6415441.5
$ws.Range("M132").Value = -6412911.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1498.5385
$ws.Range("I22").Value = 1068.4286
$ws.Range("J22").Value = 2000.3334
$ws.Range("K22").Value = 1068.4286
$ws.Range("L22").Value = 2000.3334
$ws.Range("M22").Value = -773.4286
$ws.Range("N22").Value = -2590.3334

$ws.Range("H25").Value = 9499
$ws.Range("J25").Value = 9499
$ws.Range("L25").Value = 9499
$ws.Range("N25").Value = -9959

$ws.Range("H27").Value = 1498.5385
$ws.Range("I27").Value = 1068.4286
$ws.Range("J27").Value = 2000.3334
$ws.Range("K27").Value = 1068.4286
$ws.Range("L27").Value = 2000.3334
$ws.Range("M27").Value = -961.4286
$ws.Range("N27").Value = -2214.3334

$ws.Range("H46").Value = 1127.3334
$ws.Range("I46").Value = 377
$ws.Range("J46").Value = 1377.4445
$ws.Range("K46").Value = 377
$ws.Range("L46").Value = 1377.4445
$ws.Range("M46").Value = -189
$ws.Range("N46").Value = -1753.4445

$ws.Range("H55").Value = 231.77777
$ws.Range("I55").Value = 250.58333
$ws.Range("J55").Value = 194.16667
$ws.Range("K55").Value = 250.58333
$ws.Range("L55").Value = 194.16667
$ws.Range("M55").Value = -77.58332999999999
$ws.Range("N55").Value = -540.1666700000001

$ws.Range("H100").Value = 2997.5
$ws.Range("I100").Value = 1668.6666
$ws.Range("K100").Value = 1668.6666
$ws.Range("M100").Value = -1127.6666

$ws.Range("H136").Value = 2874.6597
$ws.Range("I136").Value = 2084.8438
$ws.Range("K136").Value = 6254.5314
$ws.Range("M136").Value = -3704.5314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11339788
$ws.Range("I136").Value = 27780224
$ws.Range("K136").Value = 83340672
$ws.Range("M136").Value = -83338122
